$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a few field labels in the existing "User" / "Image" mini-tables.
$ws.Range("A7").Value = "role_id"
$ws.Range("A8").Value = "banned"
$ws.Range("C9").Value = "size_type"

# Add the new "Role" resource table (C19:C21) describing a Role entity.
$ws.Range("C19").Value = "Role"
$ws.Range("C20").Value = "id"
$ws.Range("C21").Value = "function"

$lo = $ws.ListObjects.Add(1, $ws.Range("C19:C21"), 0, 1)
$lo = $ws.ListObjects.Item("Table8")
$lo.TableStyle = "TableStyleMedium13"
$lo.Name = "Tableau92"

# Misc bookkeeping changes captured by the diff.
$null = $ws.Range("F15").Select()
$ws.PageSetup.PaperSize = 11
